$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'97.123.06"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "'3.705.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.80%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("E5").Value = "  +7.68%  "
$ws.Range("D6").Value = "'235.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.62%  "
$ws.Range("D7").Value = "'656.71"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.45%  "
$ws.Range("E8").Value = "  +1.89%  "
$ws.Range("E9").Value = "  -0.35%  "
$ws.Range("E10").Value = "  -0.03%  "
$ws.Range("D11").Value = "'3.704.63"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.86%  "
$ws.Range("D12").Value = "'45.08"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.78%  "
$ws.Range("D13").Value = "'0.0000310"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +15.36%  "
$ws.Range("D14").Value = "'0.207"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.38%  "
$ws.Range("D15").Value = "'6.89"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.16%  "
$ws.Range("D16").Value = "'4.397.96"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.82%  "
$ws.Range("D17").Value = "'97.150.98"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.42%  "
$ws.Range("E18").Value = "  +1.02%  "
$ws.Range("D19").Value = "'3.707.73"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("E20").Value = "  +1.81%  "
$ws.Range("E21").Value = "  -1.84%  "
$ws.Range("D22").Value = "'0.523"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.27%  "
$ws.Range("D23").Value = "'524.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.17%  "
$ws.Range("D24").Value = "'3.46"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.08%  "
$ws.Range("D25").Value = "'0.0000224"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +9.71%  "
$ws.Range("D26").Value = "'6.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.66%  "
$ws.Range("D27").Value = "'107.38"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.70%  "
$ws.Range("E28").Value = "  +16.63%  "
$ws.Range("D29").Value = "'3.908.59"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.94%  "
$ws.Range("D30").Value = "'13.52"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.22%  "
$ws.Range("D31").Value = "'12.64"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.32%  "
$ws.Range("E32").Value = "  -0.52%  "
$ws.Range("D33").Value = "'0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("E34").Value = "  +2.75%  "
$ws.Range("E35").Value = "  -3.48%  "
$ws.Range("E36").Value = "  +0.19%  "
$ws.Range("D37").Value = "'1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.62%  "
$ws.Range("D38").Value = "'642.75"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.97%  "
$ws.Range("D39").Value = "'0.594"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.83%  "
$ws.Range("D40").Value = "'8.73"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.87%  "
$ws.Range("E42").Value = "  +1.91%  "
$ws.Range("D43").Value = "'0.500"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +10.74%  "
$ws.Range("D44").Value = "'6.83"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.13%  "
$ws.Range("D45").Value = "'2.03"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.43%  "
$ws.Range("D46").Value = "'40.13"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.09%  "
$ws.Range("D47").Value = "'0.963"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.13%  "
$ws.Range("D48").Value = "'0.0458"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.42%  "
$ws.Range("D49").Value = "'2.40"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.96%  "
$ws.Range("E50").Value = "  -0.15%  "
$ws.Range("D51").Value = "'8.74"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.04%  "